# Bump the "Förändrad" (Changed) date in column C from 2024-08-14 (45518)
# to 2024-08-15 (45519) for every data row in the sheet (rows 2-28).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = 28

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45518) {
        $cell.Value2 = 45519
    }
}
